# Updated cryptos list on Wed Jul 31 11:58:47 UTC 2024 with GitHub Actions
#
# Applies the latest price/volume refresh to the cryptos worksheet.
# Column D ("Price") holds values that look numeric but are stored as
# literal text (e.g. "66.089.09", "0.650", "0.0000165"); we force the
# whole column to Text format before writing so Excel doesn't coerce
# them into doubles and mangle trailing zeros / thousand separators /
# scientific notation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep column D as text so values like "586.69", "0.650", "66.089.09"
# round-trip exactly instead of being reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.089.09"
$ws.Range("E2").Value = "  -0.83%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.322.28"
$ws.Range("E3").Value = "  -0.45%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "586.69"
$ws.Range("E5").Value = "  +2.26%  "

# Row 6 - Solana
$ws.Range("D6").Value = "181.70"
$ws.Range("E6").Value = "  +0.32%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.650"
$ws.Range("E7").Value = "  +2.83%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.320.72"
$ws.Range("E9").Value = "  -0.44%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.89%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +2.25%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.403"
$ws.Range("E12").Value = "  -0.57%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.898.25"
$ws.Range("E13").Value = "  -0.50%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "66.142.85"
$ws.Range("E15").Value = "  -0.90%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "26.22"
$ws.Range("E16").Value = "  -3.01%  "

# Row 17 & 18 swapped places (ShibaInu <-> WrappedEther) and got new data
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.353.13"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0000165"
$ws.Range("E18").Value = "  -1.11%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "425.02"
$ws.Range("E19").Value = "  -2.72%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  -2.60%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  -2.90%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -2.77%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "71.63"
$ws.Range("E23").Value = "  -2.47%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.01%  "

# Row 25 - LEO
$ws.Range("E25").Value = "  +0.27%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "3.461.33"
$ws.Range("E26").Value = "  -0.76%  "

# Row 27 - Polygon
$ws.Range("D27").Value = "0.513"
$ws.Range("E27").Value = "  -0.95%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +4.78%  "

# Row 29 - PEPE
$ws.Range("D29").Value = "0.0000114"
$ws.Range("E29").Value = "  -3.13%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "8.91"
$ws.Range("E30").Value = "  -1.38%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  +0.10%  "

# Row 32 - PancakeSwap
$ws.Range("D32").Value = "1.92"
$ws.Range("E32").Value = "  -1.95%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "22.39"
$ws.Range("E33").Value = "  -1.87%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  -1.94%  "

# Row 36 - Aptos
$ws.Range("E36").Value = "  -3.15%  "

# Row 37 - Fetch.AI
$ws.Range("D37").Value = "1.18"
$ws.Range("E37").Value = "  -4.04%  "

# Row 38 - Monero
$ws.Range("D38").Value = "160.82"
$ws.Range("E38").Value = "  -2.24%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -2.95%  "

# Row 40 - Maker
$ws.Range("D40").Value = "2.861.20"
$ws.Range("E40").Value = "  +0.80%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "1.80"
$ws.Range("E41").Value = "  +0.25%  "

# Row 42 - EnergySwap
$ws.Range("D42").Value = "26.41"
$ws.Range("E42").Value = "  -3.43%  "

# Row 43 - Mantle
$ws.Range("D43").Value = "0.759"
$ws.Range("E43").Value = "  -4.91%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  -2.64%  "

# Row 45 - OKB
$ws.Range("D45").Value = "39.72"
$ws.Range("E45").Value = "  -1.21%  "

# Row 46 - Hedera
$ws.Range("D46").Value = "0.0660"
$ws.Range("E46").Value = "  -1.31%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -4.92%  "

# Row 48 - dogwifhat
$ws.Range("D48").Value = "2.29"
$ws.Range("E48").Value = "  -2.20%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "23.17"
$ws.Range("E49").Value = "  -5.08%  "

# Row 50 - Bittensor
$ws.Range("D50").Value = "313.31"
$ws.Range("E50").Value = "  -2.56%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  -0.99%  "
